$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A cells from " " (space) to "x" for the requested rows
$rowsToX = @(24, 25, 26, 27, 28, 33, 34, 39, 54, 72, 73, 78, 79, 80, 81, 83, 84, 87, 101, 105, 112, 113, 114, 115, 116, 121, 122, 127, 142, 155, 160, 161, 166, 167, 168, 169, 171, 172, 175, 178, 185, 189, 190)
foreach ($r in $rowsToX) {
    $ws.Cells.Item($r, 1).Value = "x"
}

# Row 95 marker removed entirely (cell becomes blank)
$ws.Cells.Item(95, 1).Value = ""

# Reflect the saved cursor/selection position from the author (scrolled down, cell A191 active)
$ws.Activate()
$ws.Range("A185").Select()
$ws.Range("A191").Select()

Write-Host "Updated $($rowsToX.Count) cells to x and cleared A95"
